# Apply the "Exam 2 / HW TOTAL / QUIZ / TOTAL" column additions to the
# grades-m10260-sp11 roster sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Workbook window geometry (bookViews/workbookView) -------------
$wb.Windows.Item(1).Left = 2740
$wb.Windows.Item(1).Top = 7740
$wb.Windows.Item(1).Width = 21480
$wb.Windows.Item(1).Height = 4600

# --- 2. New header cells (row 1) ---------------------------------------
$ws.Range("C1").Value = "Exam 2"
$ws.Range("D1").Value = "HW TOTAL"
$ws.Range("E1").Value = "QUIZ"
$ws.Range("F1").Value = "TOTAL"

# --- 3. New data columns C (Exam 2), D (HW TOTAL), E (QUIZ) for
#        rows 2-35, bulk-assigned via a 2D array ------------------------
$arr = New-Object 'object[,]' 34,3
$arr[0,0] = 95; $arr[0,1] = 74; $arr[0,2] = 15
$arr[1,0] = 78; $arr[1,1] = 79; $arr[1,2] = 15
$arr[2,0] = 98; $arr[2,1] = 78; $arr[2,2] = 15
$arr[3,0] = 93; $arr[3,1] = 78; $arr[3,2] = 15
$arr[4,0] = 90; $arr[4,1] = 75; $arr[4,2] = 15
$arr[5,0] = 72; $arr[5,1] = 74; $arr[5,2] = 15
$arr[6,0] = 77; $arr[6,1] = 78; $arr[6,2] = 15
$arr[7,0] = 85; $arr[7,1] = 80; $arr[7,2] = 15
$arr[8,0] = 72; $arr[8,1] = 74; $arr[8,2] = 15
$arr[9,0] = 80; $arr[9,1] = 79; $arr[9,2] = 15
$arr[10,0] = 64; $arr[10,1] = 76; $arr[10,2] = 15
$arr[11,0] = 69; $arr[11,1] = 69; $arr[11,2] = 15
$arr[12,0] = 83; $arr[12,1] = 79; $arr[12,2] = 15
$arr[13,0] = 50; $arr[13,1] = 77; $arr[13,2] = 15
$arr[14,0] = 71; $arr[14,1] = 52; $arr[14,2] = 15
$arr[15,0] = 68; $arr[15,1] = 75; $arr[15,2] = 15
$arr[16,0] = 77; $arr[16,1] = 79; $arr[16,2] = 15
$arr[17,0] = 73; $arr[17,1] = 73; $arr[17,2] = 15
$arr[18,0] = 71; $arr[18,1] = 66; $arr[18,2] = 15
$arr[19,0] = 68; $arr[19,1] = 78; $arr[19,2] = 15
$arr[20,0] = 72; $arr[20,1] = 75; $arr[20,2] = 15
$arr[21,0] = 71; $arr[21,1] = 75; $arr[21,2] = 15
$arr[22,0] = 74; $arr[22,1] = 74; $arr[22,2] = 15
$arr[23,0] = 43; $arr[23,1] = 64; $arr[23,2] = 15
$arr[24,0] = 78; $arr[24,1] = 76; $arr[24,2] = 15
$arr[25,0] = 83; $arr[25,1] = 77; $arr[25,2] = 15
$arr[26,0] = 90; $arr[26,1] = 75; $arr[26,2] = 15
$arr[27,0] = 80; $arr[27,1] = 77; $arr[27,2] = 15
$arr[28,0] = 83; $arr[28,1] = 65; $arr[28,2] = 15
$arr[29,0] = 74; $arr[29,1] = 72; $arr[29,2] = 15
$arr[30,0] = 86; $arr[30,1] = 76; $arr[30,2] = 15
$arr[31,0] = 100; $arr[31,1] = 78; $arr[31,2] = 15
$arr[32,0] = 71; $arr[32,1] = 78; $arr[32,2] = 15
$arr[33,0] = 77; $arr[33,1] = 76; $arr[33,2] = 15
$ws.Range("C2:E35").Value = $arr

# --- 4. TOTAL column (F) formula. F2 is typed standalone, then F3:F35
#        filled as one block so it forms its own shared-formula group
#        (matches how the author built this up in Excel). -------------
$ws.Range("F2").Formula = "=ROUND(SUM(B2:E2)/295 * 100, 0)"
$ws.Range("F3:F35").Formula = "=ROUND(SUM(B3:E3)/295 * 100, 0)"

# --- 5. Summary rows 37 (Median) and 38 (Mean) for the new columns ----
$ws.Range("C37").Formula = "=MEDIAN(C2:C35)"
$ws.Range("D37").Formula = "=MEDIAN(D2:D35)"
$ws.Range("F37").Formula = "=MEDIAN(F2:F35)"

$ws.Range("C38").Formula = "=ROUNDUP(AVERAGE(C2:C35),1)"
$ws.Range("D38").Formula = "=ROUNDUP(AVERAGE(D2:D35),1)"
$ws.Range("F38").Formula = "=ROUNDUP(AVERAGE(F2:F35),1)"

# --- 6. View state: scroll down so row 33 is at the top, select the
#        newly-filled F37:F38 summary cells ----------------------------
$ws.Range("F37:F38").Select()

# --- 7. Page setup: portrait orientation -------------------------------
$ws.PageSetup.Orientation = 1
